$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.479.10'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').Value = '1.796.38'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = "'227.07"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('E6').Value = '  +1.67%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'32.46"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.06%  '
$ws.Range('E9').Value = '  +1.80%  '
$ws.Range('D10').Value = "'0.0694"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.87%  '
$ws.Range('D11').Value = "'0.0950"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = '2.057.36'
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('E13').Value = '  -0.40%  '
$ws.Range('D14').Value = '1.796.16'
$ws.Range('E14').Value = '  +0.70%  '
$ws.Range('D15').Value = "'0.638"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.93%  '
$ws.Range('D16').Value = '34.422.52'
$ws.Range('D17').Value = "'4.26"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').Value = "'68.74"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').Value = "'247.04"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.76%  '
$ws.Range('D20').Value = '0.0₃0801'
$ws.Range('E20').Value = '  +3.01%  '
$ws.Range('D21').Value = "'11.18"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.21%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('D24').Value = "'2.08"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.69%  '
$ws.Range('D25').Value = "'163.61"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.32%  '
$ws.Range('E26').Value = '  +1.84%  '
$ws.Range('D27').Value = "'16.56"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.66%  '
$ws.Range('E28').Value = '  +2.44%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = "'0.0523"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.04%  '
$ws.Range('E31').Value = '  +8.71%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = "'1.23"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('E33').Value = '  +3.71%  '
$ws.Range('E34').Value = '  +1.32%  '
$ws.Range('D35').Value = '1.445.96'
$ws.Range('E35').Value = '  -0.57%  '
$ws.Range('D36').Value = "'2.61"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +7.67%  '
$ws.Range('D37').Value = "'0.671"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.99%  '
$ws.Range('E38').Value = '  +2.58%  '
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('D40').Value = "'84.15"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.97%  '
$ws.Range('D41').Value = "'2.41"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.36%  '
$ws.Range('E42').Value = '  +1.73%  '
$ws.Range('E43').Value = '  +2.38%  '
$ws.Range('D44').Value = "'13.82"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.26%  '
$ws.Range('D45').Value = "'0.0526"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.36%  '
$ws.Range('D46').Value = "'6.10"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.77%  '
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').Value = '1.953.80'
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = "'105.79"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.36%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0132'
$ws.Range('E50').Value = '  -2.21%  '
$ws.Range('E51').Value = '  -0.09%  '
